$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.083.21"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.098.56"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.28"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.67"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.089.82"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("E11").Value = "  +6.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.10"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").Value = "3.590.82"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "63.083.35"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "3.089.94"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "503.89"
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.68"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.52"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.28"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.58"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.35"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.26"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.30"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.52"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.12"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "59.26"
$ws.Range("E34").Value = "  +14.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "528.99"
$ws.Range("E35").Value = "  -7.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  +2.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.22"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0415"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0791"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.121"
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("D41").Value = "3.040.17"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.67"
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.08"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.254"
$ws.Range("E44").Value = "  +6.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.31"
$ws.Range("E47").Value = "  +4.31%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.94"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.106"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.44"
$ws.Range("E50").Value = "  +75.52%  "
$ws.Range("D51").Value = "0.0₃0503"
$ws.Range("E51").Value = "  -1.95%  "
